$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows(24).Insert()
$ws.Range("A24").Value = "get distance between locations"
$ws.Range("B24:C24").Value = "NOT STARTED"
$ws.Range("B25:C25").Copy()
$ws.Range("B24:C24").PasteSpecial(-4122)
$ws.Rows("41:42").Insert()
